$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data area first
$ws.Range("A1:C4").Clear()

# New header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "username_value"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "password_value"
$ws.Range("E1").Value = "login_button"
$ws.Range("F1").Value = "exp_login"

# New data row
$ws.Range("A2").Value = "//input[@id ='user-name']"
$ws.Range("B2").Value = "standard_user"
$ws.Range("C2").Value = "//input[@id ='password']"
$ws.Range("D2").Value = "secret_sauce"
$ws.Range("E2").Value = "//input[@id ='login-button']"
$ws.Range("F2").Value = "Yes"

# Select F2 as the active selection (matches diff's <selection activeCell="F2" sqref="F2"/>)
$ws.Range("F2").Select()

# Auto-fit columns to match bestFit widths in target (only A:E get custom widths)
$ws.Columns("A:E").AutoFit()

# Page setup: portrait orientation (matches new <pageSetup orientation="portrait".../>)
$ws.PageSetup.Orientation = 1
